$d = $word.ActiveDocument

# Find the "Docente(s) Responsável(eis) " heading paragraph
$findRange = $d.Content
$found = $findRange.Find.Execute("Docente(s) Responsável(eis) ", $false, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$headingPara = $findRange.Paragraphs(1)

# Insert a new (blank, inherits heading style) paragraph right after the heading
$headingPara.Range.InsertParagraphAfter() | Out-Null

# The paragraph that now immediately follows the heading is the freshly inserted one
$newPara = $headingPara.Next()

# Build the new paragraph body as raw OOXML so the two names land in separate
# runs joined by a line break, matching the target structure exactly:
#   <w:r><w:t>5817330 - Larissa de Freitas</w:t><w:br/></w:r>
#   <w:r><w:t>1506103 - Pedro Carlos de Oliveira</w:t></w:r>
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>5817330 - Larissa de Freitas</w:t><w:br/></w:r>' +
       '<w:r><w:t>1506103 - Pedro Carlos de Oliveira</w:t></w:r>' +
       '</w:p>'

$newPara.Range.InsertXML($xml) | Out-Null
